$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-499).
# All rows are updated from serial 45190 (2023-09-21) to 45192 (2023-09-23).
for ($row = 2; $row -le 499; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
